$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Constants used for alignment
$xlCenter = -4108

# ---------------------------------------------------------------------------
# Row 15: "saldo" header above the new table, bold + centered (C15)
# ---------------------------------------------------------------------------
$ws.Range("C15").HorizontalAlignment = $xlCenter
$ws.Range("C15").Font.Bold = $true
$ws.Range("C15").Value = "saldo"

# ---------------------------------------------------------------------------
# Row 16: year headers (2016..2022) built with relative formulas, then the
# summary column headers (sum / sum3 / alltime % / last 3yr %)
# ---------------------------------------------------------------------------
$ws.Range("C16:M16").HorizontalAlignment = $xlCenter
$ws.Range("C16:M16").Font.Bold = $true

$ws.Range("C16").Value = 2016
$ws.Range("D16:I16").Formula = "=C16+1"

$ws.Range("J16").Value = "sum"
$ws.Range("K16").Value = "sum3"
$ws.Range("L16").Value = "alltime %"
$ws.Range("M16").Value = "last 3yr %"

# ---------------------------------------------------------------------------
# Rows 17-25: "total LO" + each district's net migration (saldo) per year,
# pulled from the inflow/outflow block above (rows 4-12), plus sums and
# shares of the "total LO" row.
# ---------------------------------------------------------------------------
# label (col B), source row (inflow/outflow pair row) for each target row
$rows = @(
    @{ Row = 17; Label = "total LO";          Src = 4  },
    @{ Row = 18; Label = "Central dist.";      Src = 5  },
    @{ Row = 19; Label = "North-west dist.";   Src = 6  },
    @{ Row = 20; Label = "South dist.";        Src = 7  },
    @{ Row = 21; Label = "North-Caucas dist."; Src = 8  },
    @{ Row = 22; Label = "Volga dist.";        Src = 9  },
    @{ Row = 23; Label = "Ural dist.";         Src = 10 },
    @{ Row = 24; Label = "Siberia dist.";      Src = 11 },
    @{ Row = 25; Label = "Far-east dist.";     Src = 12 }
)

$targetCols = @("C","D","E","F","G","H","I")
$srcCols    = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P")

foreach ($r in $rows) {
    $row = $r.Row
    $src = $r.Src

    $ws.Range("B$row").HorizontalAlignment = $xlCenter
    $ws.Range("B$row").Font.Bold = $true
    $ws.Range("B$row").Value = $r.Label

    for ($i = 0; $i -lt $targetCols.Length; $i++) {
        $tCol = $targetCols[$i]
        $lCol = $srcCols[2 * $i]
        $rCol = $srcCols[2 * $i + 1]
        $ws.Range("$tCol$row").Formula = "=$lCol$src-$rCol$src"
    }

    $ws.Range("J$row").HorizontalAlignment = $xlCenter
    $ws.Range("J$row").Formula = "=SUM(C$row`:I$row)"

    $ws.Range("K$row").HorizontalAlignment = $xlCenter
    $ws.Range("K$row").Formula = "=SUM(G$row`:I$row)"

    if ($row -ne 17) {
        $ws.Range("L$row").HorizontalAlignment = $xlCenter
        $ws.Range("L$row").Font.Bold = $true
        $ws.Range("L$row").NumberFormat = "0.00"
        $ws.Range("L$row").Formula = "=J$row/J`$17"

        $ws.Range("M$row").HorizontalAlignment = $xlCenter
        $ws.Range("M$row").Font.Bold = $true
        $ws.Range("M$row").NumberFormat = "0.00"
        $ws.Range("M$row").Formula = "=K$row/K`$17"
    }
}

# ---------------------------------------------------------------------------
# Row 26: a lone formatted blank spacer cell (mirrors the sheet's earlier
# blank spacer rows) so the used range grows to row 26.
# ---------------------------------------------------------------------------
$ws.Range("C26").HorizontalAlignment = $ws.Range("C13").HorizontalAlignment

# ---------------------------------------------------------------------------
# Column M needs to be a bit wider to fit "last 3yr %"
# ---------------------------------------------------------------------------
$ws.Columns("M").ColumnWidth = 10.14

# ---------------------------------------------------------------------------
# Leave the selection where the author's session ended up
# ---------------------------------------------------------------------------
$null = $ws.Range("V22").Select()
